$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The alcohol measurement sheet had an extra data column (old column M)
# which duplicated/obsoleted the data that used to live in column N.
# Removing column M shifts the former column N left into its place,
# matching the new A1:M119 data range.
$ws.Columns("M").Delete()

# Reflect the resulting selection on the now-last column of data.
[void]$ws.Range("M1").Select()
